# Item.xlsx update:
#  - remove the two "Consumable" (TYPE=2) attack items: ID 1204 (화염병 / Molotov
#    cocktail, row 14) and ID 2200 (짱돌 / rock, row 22)
#  - change the Consumable "recover" items' (ID 1200-1203, rows 10-13) status
#    effect: VALUE (price) raised and AMOUNT now carries the heal amount
#    instead of a flat 1; the now-unused SKILL id (column H) is cleared
#  - every other item's AMOUNT (column E) becomes 0 instead of 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the removed attack-item rows bottom-up so earlier row numbers stay valid.
$ws.Rows(22).Delete()   # ID 2200 "짱돌"
$ws.Rows(14).Delete()   # ID 1204 "화염병"

# After the two deletions the data now spans rows 2-20.
$ws.Range("E2:E20").Value2 = 0

# Recover potions (ID 1200-1203 / rows 10-13): new VALUE + AMOUNT, drop SKILL id.
$ws.Range("D10").Value2 = 100
$ws.Range("E10").Value2 = 45
$ws.Range("H10").Clear()

$ws.Range("D11").Value2 = 250
$ws.Range("E11").Value2 = 70
$ws.Range("H11").Clear()

$ws.Range("D12").Value2 = 600
$ws.Range("E12").Value2 = 100
$ws.Range("H12").Clear()

$ws.Range("D13").Value2 = 1500
$ws.Range("E13").Value2 = 150
$ws.Range("H13").Clear()

# Match the saved view's selection.
$ws.Range("F14").Select()
